# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2404   (the "before" / FV2404 format version columns)
#   *_new -> *_FV2410   (the "after"  / FV2410 format version columns)
# Then wrap the data range in a proper Excel table (Table1) and freeze
# the header row, matching the regenerated AHB-diff export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) -----------------------------------
# Columns A:J  -> *_old  => *_FV2404
# Column  K    -> "diff" (unchanged)
# Columns L:U  -> *_new  => *_FV2410
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Text
    if ($text -like "*_old") {
        $cell.Value = ($text -replace "_old$", "_FV2404")
    } elseif ($text -like "*_new") {
        $cell.Value = ($text -replace "_new$", "_FV2410")
    }
}

# --- 2. Turn the data range into an Excel table ------------------------
$dataRange = $ws.Range("A1:U79")
$lo = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
